$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: copy formats for newly-populated / restyled rows (based on ORIGINAL layout) ----
# Order matters: copy each destination before its source row is overwritten below.

# new row 9 <- original row 8 (style 7, uniform A:G)
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)

# new row 8 (A:F) <- original row 7 (style 5, uniform A:F)
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)

# row 7 <- original row 6 (style 4, uniform A:G) -- row 7 changes style 5 -> 4
$ws.Range("A6:G6").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)

# new row 5 <- original row 4 (style 3, uniform A:G)
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# G8 needs the purple "exception" fill (same fill color G7 originally had)
$ws.Range("G8").Interior.Color = 10498160

# ---- Step 2: column H ----
$ws.Columns.Item(8).ColumnWidth = 21.166666666666668
# give each H cell (rows 2-9) the same format as its row (copy from column G)
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("G7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Step 3: write cell values, columns A-G first (top to bottom) so the shared-string
#      table gets new entries appended in the same order as the target workbook ----
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "preparationTime"
$ws.Range("D1").Value = "servings"
$ws.Range("E1").Value = "categories"
$ws.Range("F1").Value = "comments"
$ws.Range("G1").Value = "image"
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "len > 0"
$ws.Range("C2").Value = "> 0"
$ws.Range("D2").Value = "> 0"
$ws.Range("E2").Value = "len > 0"
$ws.Range("F2").Value = "X"
$ws.Range("G2").Value = "null"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "len = 0"
$ws.Range("C3").Value = "> 0"
$ws.Range("D3").Value = "> 0"
$ws.Range("E3").Value = "len > 0"
$ws.Range("F3").Value = "X"
$ws.Range("G3").Value = "null"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "len > 0"
$ws.Range("C4").Value = "< 0"
$ws.Range("D4").Value = "> 0"
$ws.Range("E4").Value = "len > 0"
$ws.Range("F4").Value = "X"
$ws.Range("G4").Value = "null"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "len > 0"
$ws.Range("C5").Value = "equals(0)"
$ws.Range("D5").Value = "> 0"
$ws.Range("E5").Value = "len > 0"
$ws.Range("F5").Value = "X"
$ws.Range("G5").Value = "null"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "len > 0"
$ws.Range("C6").Value = "> 0"
$ws.Range("D6").Value = "equals(0)"
$ws.Range("E6").Value = "len > 0"
$ws.Range("F6").Value = "X"
$ws.Range("G6").Value = "null"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "len > 0"
$ws.Range("C7").Value = "> 0"
$ws.Range("D7").Value = "< 0"
$ws.Range("E7").Value = "len > 0"
$ws.Range("F7").Value = "X"
$ws.Range("G7").Value = "null"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "len > 0"
$ws.Range("C8").Value = "> 0"
$ws.Range("D8").Value = "> 0"
$ws.Range("E8").Value = "len = 0"
$ws.Range("F8").Value = "X"
$ws.Range("G8").Value = "null"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "len > 0"
$ws.Range("C9").Value = "> 0"
$ws.Range("D9").Value = "> 0"
$ws.Range("E9").Value = "len > 0"
$ws.Range("F9").Value = "X"
$ws.Range("G9").Value = "image"

# ---- Step 4: write column H (Expected Result) top to bottom ----
$ws.Range("H1").Value = "Expected Result"
$ws.Range("H2").Value = "No Action"
$ws.Range("H3").Value = "IllegalArgumentException"
$ws.Range("H4").Value = "IllegalArgumentException"
$ws.Range("H5").Value = "IllegalArgumentException"
$ws.Range("H6").Value = "IllegalArgumentException"
$ws.Range("H7").Value = "IllegalArgumentException"
$ws.Range("H8").Value = "IllegalArgumentException"
$ws.Range("H9").Value = "No Action"

# ---- Step 5: second table (Field / Number of States) updates ----
$ws.Range("B14").Value = 3
$ws.Range("B15").Value = 3

# ---- Step 6: selection (matches final saved state) ----
$ws.Range("I4").Select()
